$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E).
# Numeric-looking Price values must be written as literal TEXT (matching the
# source data, which stores every Price/Volume cell as a string) rather than
# being auto-converted to a number by Excel. We force text entry by switching
# the cell to the "@" (Text) number format before assigning the value, then
# call ClearFormats() to drop the now-unneeded style again so the cell keeps
# its original (default) style/appearance.

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "42.041.67"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "2.299.94"
$ws.Range("E3").Value = "  -3.02%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "310.88"
$ws.Range("E5").Value = "  -7.00%  "
Set-TextValue "D6" "105.16"
$ws.Range("E6").Value = "  +4.32%  "
Set-TextValue "D7" "0.625"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -5.08%  "
Set-TextValue "D10" "39.79"
$ws.Range("E10").Value = "  -1.02%  "
Set-TextValue "D11" "0.0910"
$ws.Range("E11").Value = "  -1.26%  "
Set-TextValue "D12" "8.27"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("E13").Value = "  -0.21%  "
Set-TextValue "D14" "0.966"
$ws.Range("E14").Value = "  -4.43%  "
$ws.Range("E15").Value = "  -6.02%  "
$ws.Range("D16").Value = "2.648.89"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "2.303.17"
$ws.Range("E17").Value = "  -3.16%  "
$ws.Range("D18").Value = "41.995.40"
$ws.Range("E18").Value = "  -1.69%  "
Set-TextValue "D19" "7.59"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("E20").Value = "  -2.56%  "
Set-TextValue "D21" "74.25"
$ws.Range("E21").Value = "  -1.96%  "
Set-TextValue "D22" "3.47"
$ws.Range("E22").Value = "  -9.28%  "
Set-TextValue "D23" "258.46"
$ws.Range("E23").Value = "  -4.91%  "
$ws.Range("E24").Value = "  -3.01%  "
Set-TextValue "D25" "9.18"
$ws.Range("E25").Value = "  -7.67%  "
$ws.Range("E26").Value = "  +0.52%  "
Set-TextValue "D27" "10.96"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("E28").Value = "  +2.62%  "
Set-TextValue "D29" "22.77"
$ws.Range("E29").Value = "  -2.60%  "
Set-TextValue "D30" "35.68"
$ws.Range("E30").Value = "  +0.28%  "
Set-TextValue "D31" "163.73"
$ws.Range("E31").Value = "  -6.61%  "
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("E33").Value = "  -6.03%  "
Set-TextValue "D34" "5.83"
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("E35").Value = "  -4.25%  "
$ws.Range("E36").Value = "  +9.33%  "
$ws.Range("E37").Value = "  -2.86%  "
Set-TextValue "D38" "0.0350"
$ws.Range("E38").Value = "  -2.85%  "
Set-TextValue "D39" "3.65"
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("E40").Value = "  -7.49%  "
Set-TextValue "D41" "98.57"
$ws.Range("E41").Value = "  +9.44%  "
$ws.Range("E42").Value = "  -4.75%  "
Set-TextValue "D43" "69.63"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("E45").Value = "  -0.05%  "
Set-TextValue "D46" "12.10"
$ws.Range("E46").Value = "  -0.05%  "
Set-TextValue "D47" "110.49"
$ws.Range("E47").Value = "  -6.76%  "
Set-TextValue "D48" "5.35"
$ws.Range("E48").Value = "  -2.81%  "
Set-TextValue "D49" "8.92"
$ws.Range("E49").Value = "  -2.58%  "
Set-TextValue "D50" "72.86"
$ws.Range("E50").Value = "  +3.42%  "
Set-TextValue "D51" "1.26"
$ws.Range("E51").Value = "  -1.45%  "
